# Update workbook from "2021-01-01" edition to "2021-01-02" edition.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the data sheet (also updates the 'Bundesländer001' defined
#    name, which refers to this sheet by name).
# ------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("31.12.20")
$wsData.Name = "01.01.21"

$wsInfo = $wb.Worksheets.Item("Erläuterung")

# ------------------------------------------------------------------
# 2. "Erläuterung" sheet text updates.
# ------------------------------------------------------------------
# Time-of-day changes from 12:30 Uhr to 08:00 Uhr.
$wsInfo.Range("C6").Value = "08:00 Uhr"

# Footnote text drops "oder Korrekturen".
$wsInfo.Range("A10").Value = "Achtung: Die Differenz zum Vortag kann Nachmeldungen aus vorangegangenen Tagen enthalten und spiegelt nicht immer die innerhalb des Vortags tatsächlich durchgeführte Zahl der Impfungen wider. "

# ------------------------------------------------------------------
# 3. Data sheet ("01.01.21") value updates.
# ------------------------------------------------------------------

# Header indication columns keep the same text (shared-string reindex
# only in the source diff) - nothing to change here.

# Row 2: Baden-Württemberg
$wsData.Range("B2").Value = 20045
$wsData.Range("C2").Value = 2958
$wsData.Range("D2").Value = 9087
$wsData.Range("E2").Value = 6185
$wsData.Range("F2").Value = 1189
$wsData.Range("G2").Value = 3403

# Row 3: Bayern
$wsData.Range("B3").Value = 39005
$wsData.Range("C3").Value = 1050
$wsData.Range("E3").Value = 19500
$wsData.Range("H3").Value = "(Nachmeldungen für den 01.01. bereits für morgen angekündigt)"

# Row 4: Berlin
$wsData.Range("B4").Value = 13137
$wsData.Range("C4").Value = 2023
$wsData.Range("D4").Value = 8662
$wsData.Range("E4").Value = 3319
$wsData.Range("G4").Value = 9817

# Row 5: Brandenburg
$wsData.Range("C5").Value = ""
$wsData.Range("H5").Value = "(für den 01.01. wurde nicht gemeldet) "

# Row 6: Bremen
$wsData.Range("C6").Value = ""

# Row 7: Hamburg
$wsData.Range("B7").Value = 3042
$wsData.Range("C7").Value = 283
$wsData.Range("E7").Value = 1703

# Row 8: Hessen
$wsData.Range("B8").Value = 24791
$wsData.Range("C8").Value = 3418
$wsData.Range("D8").Value = 6909
$wsData.Range("E8").Value = 14480
$wsData.Range("F8").Value = 955
$wsData.Range("G8").Value = 9680

# Row 9: Mecklenburg-Vorpommern
$wsData.Range("C9").Value = ""

# Row 10: Niedersachsen
$wsData.Range("B10").Value = 3945
$wsData.Range("C10").Value = 379
$wsData.Range("D10").Value = 773
$wsData.Range("E10").Value = 2105
$wsData.Range("F10").Value = 827
$wsData.Range("G10").Value = 2362

# Row 11: Nordrhein-Westfalen
$wsData.Range("B11").Value = 33375
$wsData.Range("C11").Value = 7908
$wsData.Range("E11").Value = 11846
$wsData.Range("F11").Value = 0
$wsData.Range("G11").Value = 21537
$wsData.Range("H11").Value = ""

# Row 12: Rheinland-Pfalz
$wsData.Range("B12").Value = 6898
$wsData.Range("C12").Value = 439
$wsData.Range("E12").Value = 3166
$wsData.Range("G12").Value = 3732

# Row 13: Saarland
$wsData.Range("B13").Value = 3316
$wsData.Range("C13").Value = 600
$wsData.Range("D13").Value = 2508
$wsData.Range("E13").Value = 335
$wsData.Range("G13").Value = 1294

# Row 14: Sachsen
$wsData.Range("B14").Value = 4000
$wsData.Range("C14").Value = 710
$wsData.Range("D14").Value = 297
$wsData.Range("E14").Value = 3320
$wsData.Range("G14").Value = 679

# Row 15: Sachsen-Anhalt
$wsData.Range("B15").Value = 11771
$wsData.Range("C15").Value = 625
$wsData.Range("D15").Value = 3782
$wsData.Range("E15").Value = 5760
$wsData.Range("F15").Value = 495
$wsData.Range("G15").Value = 5793

# Row 16: Schleswig-Holstein
$wsData.Range("B16").Value = 7964
$wsData.Range("C16").Value = 694
$wsData.Range("D16").Value = 2130
$wsData.Range("E16").Value = 4690
$wsData.Range("F16").Value = 1798
$wsData.Range("G16").Value = 3384

# Row 17: Thüringen
$wsData.Range("C17").Value = ""
$wsData.Range("F17").Value = ""

# Row 18 contains SUM() formulas already and will recalc automatically.

# ------------------------------------------------------------------
# 4. Restore the active-cell selections that Excel records per sheet.
# ------------------------------------------------------------------
[void]$wsInfo.Range("C6").Select()
[void]$wsData.Range("H3").Select()
